$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1498
$ws.Range("K3").Value = 1434
$ws.Range("K4").Value = 308
$ws.Range("K6").Value = 1834
$ws.Range("K7").Value = 5170

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 99
$ws.Range("K3").Value = 94
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 322

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J4").Value = 22
$ws.Range("J7").Value = 590

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 205

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 58
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 37
$ws.Range("K8").Value = 322
$ws.Range("K11").Value = 105
$ws.Range("K19").Value = 133
$ws.Range("K20").Value = 119
$ws.Range("K23").Value = 53
$ws.Range("K29").Value = 234
$ws.Range("K31").Value = 61
$ws.Range("K33").Value = 205
$ws.Range("K37").Value = 176
$ws.Range("K42").Value = 180
$ws.Range("K43").Value = 50
$ws.Range("K47").Value = 38
$ws.Range("K51").Value = 60
$ws.Range("K52").Value = 140
$ws.Range("K54").Value = 89
$ws.Range("K55").Value = 55
$ws.Range("K60").Value = 40
$ws.Range("K61").Value = 9
$ws.Range("J63").Value = 90
$ws.Range("K63").Value = 17
$ws.Range("K67").Value = 199
$ws.Range("K72").Value = 23
$ws.Range("K77").Value = 38
$ws.Range("K79").Value = 138
$ws.Range("J83").Value = 590
$ws.Range("K84").Value = 35
$ws.Range("J85").Value = 1193
$ws.Range("K85").Value = 270
$ws.Range("J86").Value = 176
$ws.Range("K87").Value = 4
$ws.Range("K88").Value = 66
$ws.Range("K89").Value = 67
$ws.Range("K91").Value = 51
$ws.Range("K93").Value = 22
$ws.Range("K94").Value = 63
$ws.Range("K95").Value = 89
$ws.Range("K96").Value = 68
$ws.Range("K97").Value = 44
$ws.Range("K99").Value = 95
$ws.Range("K101").Value = 5170

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 9
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 57
$ws.Range("K3").Value = 62
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 31
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 77
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 234

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 43
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 42
$ws.Range("K3").Value = 48
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 11
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 35
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 33
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 9
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 8
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 95
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 98
$ws.Range("K3").Value = 85
$ws.Range("J4").Value = 72
$ws.Range("K4").Value = 15
$ws.Range("K6").Value = 68
$ws.Range("J7").Value = 1193
$ws.Range("K7").Value = 270

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 4

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 9
